# Auto-generated edit script for primera-rfef-group-2 2023-2024 workbook
# Reorders betting-odds rows that share a kickoff slot back into their
# correct fixture order, and appends the newly scraped Sanluqueno vs
# Cordoba fixture as row 128.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: restore correct fixture data (F:V)
$ws.Cells.Item(18, 6).Value = "Antequera"
$ws.Cells.Item(18, 7).Value = 1
$ws.Cells.Item(18, 8).Value = "Real Madrid B"
$ws.Cells.Item(18, 9).Value = 2
$ws.Cells.Item(18, 10).Value = 2.63
$ws.Cells.Item(18, 11).Value = "01/09/2023 17:43"
$ws.Cells.Item(18, 12).Value = 2.86
$ws.Cells.Item(18, 13).Value = "03/09/2023 21:27"
$ws.Cells.Item(18, 14).Value = 3.1
$ws.Cells.Item(18, 15).Value = "01/09/2023 17:43"
$ws.Cells.Item(18, 16).Value = 3.1
$ws.Cells.Item(18, 17).Value = "03/09/2023 21:27"
$ws.Cells.Item(18, 18).Value = 2.51
$ws.Cells.Item(18, 19).Value = "01/09/2023 17:43"
$ws.Cells.Item(18, 20).Value = 2.6
$ws.Cells.Item(18, 21).Value = "03/09/2023 21:27"
$ws.Cells.Item(18, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/antequera-real-madrid/YyqloyjL/"

# Row 20: restore correct fixture data (F:V)
$ws.Cells.Item(20, 6).Value = "CF Intercity"
$ws.Cells.Item(20, 7).Value = 1
$ws.Cells.Item(20, 8).Value = "Alcoyano"
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = 2.14
$ws.Cells.Item(20, 11).Value = "01/09/2023 17:43"
$ws.Cells.Item(20, 12).Value = 2.52
$ws.Cells.Item(20, 13).Value = "03/09/2023 19:01"
$ws.Cells.Item(20, 14).Value = 2.97
$ws.Cells.Item(20, 15).Value = "01/09/2023 17:43"
$ws.Cells.Item(20, 16).Value = 3.0
$ws.Cells.Item(20, 17).Value = "03/09/2023 19:31"
$ws.Cells.Item(20, 18).Value = 3.36
$ws.Cells.Item(20, 19).Value = "01/09/2023 17:43"
$ws.Cells.Item(20, 20).Value = 3.02
$ws.Cells.Item(20, 21).Value = "03/09/2023 19:01"
$ws.Cells.Item(20, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/cf-intercity-alcoyano/nqPOvZSl/"

# Row 65: restore correct fixture data (F:V)
$ws.Cells.Item(65, 6).Value = "Baleares"
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = "Murcia"
$ws.Cells.Item(65, 9).Value = 1
$ws.Cells.Item(65, 10).Value = 3.21
$ws.Cells.Item(65, 11).Value = "05/10/2023 08:13"
$ws.Cells.Item(65, 12).Value = 3.42
$ws.Cells.Item(65, 13).Value = "07/10/2023 13:28"
$ws.Cells.Item(65, 14).Value = 2.94
$ws.Cells.Item(65, 15).Value = "05/10/2023 08:13"
$ws.Cells.Item(65, 16).Value = 3.04
$ws.Cells.Item(65, 17).Value = "08/10/2023 10:02"
$ws.Cells.Item(65, 18).Value = 2.23
$ws.Cells.Item(65, 19).Value = "05/10/2023 08:13"
$ws.Cells.Item(65, 20).Value = 2.24
$ws.Cells.Item(65, 21).Value = "08/10/2023 06:53"
$ws.Cells.Item(65, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/baleares-murcia/8Y7GOGwp/"

# Row 66: restore correct fixture data (F:V)
$ws.Cells.Item(66, 6).Value = "Linares"
$ws.Cells.Item(66, 7).Value = 2
$ws.Cells.Item(66, 8).Value = "Ceuta"
$ws.Cells.Item(66, 9).Value = 1
$ws.Cells.Item(66, 10).Value = 3.5
$ws.Cells.Item(66, 11).Value = "05/10/2023 08:13"
$ws.Cells.Item(66, 12).Value = 3.19
$ws.Cells.Item(66, 13).Value = "08/10/2023 11:26"
$ws.Cells.Item(66, 14).Value = 3.05
$ws.Cells.Item(66, 15).Value = "05/10/2023 08:13"
$ws.Cells.Item(66, 16).Value = 3.15
$ws.Cells.Item(66, 17).Value = "08/10/2023 10:51"
$ws.Cells.Item(66, 18).Value = 2.09
$ws.Cells.Item(66, 19).Value = "05/10/2023 08:13"
$ws.Cells.Item(66, 20).Value = 2.34
$ws.Cells.Item(66, 21).Value = "08/10/2023 11:26"
$ws.Cells.Item(66, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/linares-ad-ceuta/jHmFqIF9/"

# Row 75: restore correct fixture data (F:V)
$ws.Cells.Item(75, 6).Value = "Ceuta"
$ws.Cells.Item(75, 7).Value = 1
$ws.Cells.Item(75, 8).Value = "Castellon"
$ws.Cells.Item(75, 9).Value = 3
$ws.Cells.Item(75, 10).Value = 2.9
$ws.Cells.Item(75, 11).Value = "12/10/2023 08:13"
$ws.Cells.Item(75, 12).Value = 2.92
$ws.Cells.Item(75, 13).Value = "15/10/2023 11:54"
$ws.Cells.Item(75, 14).Value = 2.96
$ws.Cells.Item(75, 15).Value = "12/10/2023 08:13"
$ws.Cells.Item(75, 16).Value = 3.3
$ws.Cells.Item(75, 17).Value = "15/10/2023 11:54"
$ws.Cells.Item(75, 18).Value = 2.45
$ws.Cells.Item(75, 19).Value = "12/10/2023 08:13"
$ws.Cells.Item(75, 20).Value = 2.43
$ws.Cells.Item(75, 21).Value = "15/10/2023 11:54"
$ws.Cells.Item(75, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/ad-ceuta-castellon/xEL2zExd/"

# Row 76: restore correct fixture data (F:V)
$ws.Cells.Item(76, 6).Value = "Murcia"
$ws.Cells.Item(76, 7).Value = 1
$ws.Cells.Item(76, 8).Value = "Algeciras"
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 10).Value = 1.97
$ws.Cells.Item(76, 11).Value = "12/10/2023 08:13"
$ws.Cells.Item(76, 12).Value = 2.37
$ws.Cells.Item(76, 13).Value = "15/10/2023 11:58"
$ws.Cells.Item(76, 14).Value = 3.07
$ws.Cells.Item(76, 15).Value = "12/10/2023 08:13"
$ws.Cells.Item(76, 16).Value = 2.98
$ws.Cells.Item(76, 17).Value = "15/10/2023 11:58"
$ws.Cells.Item(76, 18).Value = 3.74
$ws.Cells.Item(76, 19).Value = "12/10/2023 08:13"
$ws.Cells.Item(76, 20).Value = 3.33
$ws.Cells.Item(76, 21).Value = "15/10/2023 11:58"
$ws.Cells.Item(76, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/murcia-algeciras/4KDlbHU2/"

# Row 77: restore correct fixture data (F:V)
$ws.Cells.Item(77, 6).Value = "UD Ibiza"
$ws.Cells.Item(77, 7).Value = 5
$ws.Cells.Item(77, 8).Value = "Recreativo Huelva"
$ws.Cells.Item(77, 9).Value = 2
$ws.Cells.Item(77, 10).Value = 1.59
$ws.Cells.Item(77, 11).Value = "12/10/2023 08:13"
$ws.Cells.Item(77, 12).Value = 1.54
$ws.Cells.Item(77, 13).Value = "15/10/2023 11:54"
$ws.Cells.Item(77, 14).Value = 3.55
$ws.Cells.Item(77, 15).Value = "12/10/2023 08:13"
$ws.Cells.Item(77, 16).Value = 3.85
$ws.Cells.Item(77, 17).Value = "15/10/2023 11:54"
$ws.Cells.Item(77, 18).Value = 5.47
$ws.Cells.Item(77, 19).Value = "12/10/2023 08:13"
$ws.Cells.Item(77, 20).Value = 6.68
$ws.Cells.Item(77, 21).Value = "15/10/2023 11:54"
$ws.Cells.Item(77, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/ud-ibiza-recreativo-huelva/AXrkdeaF/"

# Row 84: restore correct fixture data (F:V)
$ws.Cells.Item(84, 6).Value = "CF Intercity"
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = "Baleares"
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 1.75
$ws.Cells.Item(84, 11).Value = "19/10/2023 08:13"
$ws.Cells.Item(84, 12).Value = 1.75
$ws.Cells.Item(84, 13).Value = "22/10/2023 09:56"
$ws.Cells.Item(84, 14).Value = 3.24
$ws.Cells.Item(84, 15).Value = "19/10/2023 08:13"
$ws.Cells.Item(84, 16).Value = 3.42
$ws.Cells.Item(84, 17).Value = "22/10/2023 10:01"
$ws.Cells.Item(84, 18).Value = 4.53
$ws.Cells.Item(84, 19).Value = "19/10/2023 08:13"
$ws.Cells.Item(84, 20).Value = 5.02
$ws.Cells.Item(84, 21).Value = "22/10/2023 09:56"
$ws.Cells.Item(84, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/cf-intercity-baleares/djIgxePt/"

# Row 85: restore correct fixture data (F:V)
$ws.Cells.Item(85, 6).Value = "Granada CF B"
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = "UD Ibiza"
$ws.Cells.Item(85, 9).Value = 1
$ws.Cells.Item(85, 10).Value = 4.18
$ws.Cells.Item(85, 11).Value = "19/10/2023 08:13"
$ws.Cells.Item(85, 12).Value = 4.16
$ws.Cells.Item(85, 13).Value = "22/10/2023 10:38"
$ws.Cells.Item(85, 14).Value = 3.2
$ws.Cells.Item(85, 15).Value = "19/10/2023 08:13"
$ws.Cells.Item(85, 16).Value = 3.14
$ws.Cells.Item(85, 17).Value = "22/10/2023 10:03"
$ws.Cells.Item(85, 18).Value = 1.85
$ws.Cells.Item(85, 19).Value = "19/10/2023 08:13"
$ws.Cells.Item(85, 20).Value = 2.0
$ws.Cells.Item(85, 21).Value = "22/10/2023 04:23"
$ws.Cells.Item(85, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/granada-cf-ud-ibiza/tK5Mwmql/"

# Row 86: restore correct fixture data (F:V)
$ws.Cells.Item(86, 6).Value = "Melilla"
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = "Antequera"
$ws.Cells.Item(86, 9).Value = 1
$ws.Cells.Item(86, 10).Value = 2.96
$ws.Cells.Item(86, 11).Value = "19/10/2023 08:13"
$ws.Cells.Item(86, 12).Value = 3.0
$ws.Cells.Item(86, 13).Value = "22/10/2023 06:00"
$ws.Cells.Item(86, 14).Value = 2.91
$ws.Cells.Item(86, 15).Value = "19/10/2023 08:13"
$ws.Cells.Item(86, 16).Value = 2.91
$ws.Cells.Item(86, 17).Value = "22/10/2023 10:01"
$ws.Cells.Item(86, 18).Value = 2.44
$ws.Cells.Item(86, 19).Value = "19/10/2023 08:13"
$ws.Cells.Item(86, 20).Value = 2.59
$ws.Cells.Item(86, 21).Value = "22/10/2023 06:00"
$ws.Cells.Item(86, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/melilla-antequera/vBVRtyHP/"

# Row 101: restore correct fixture data (F:V)
$ws.Cells.Item(101, 6).Value = "Linares"
$ws.Cells.Item(101, 7).Value = 1
$ws.Cells.Item(101, 8).Value = "Algeciras"
$ws.Cells.Item(101, 9).Value = 1
$ws.Cells.Item(101, 10).Value = 2.75
$ws.Cells.Item(101, 11).Value = "02/11/2023 08:13"
$ws.Cells.Item(101, 12).Value = 3.01
$ws.Cells.Item(101, 13).Value = "04/11/2023 19:52"
$ws.Cells.Item(101, 14).Value = 2.8
$ws.Cells.Item(101, 15).Value = "02/11/2023 08:13"
$ws.Cells.Item(101, 16).Value = 2.98
$ws.Cells.Item(101, 17).Value = "04/11/2023 19:52"
$ws.Cells.Item(101, 18).Value = 2.63
$ws.Cells.Item(101, 19).Value = "02/11/2023 08:13"
$ws.Cells.Item(101, 20).Value = 2.56
$ws.Cells.Item(101, 21).Value = "04/11/2023 19:52"
$ws.Cells.Item(101, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/linares-algeciras/nJu4rlhU/"

# Row 102: restore correct fixture data (F:V)
$ws.Cells.Item(102, 6).Value = "CF Intercity"
$ws.Cells.Item(102, 7).Value = 1
$ws.Cells.Item(102, 8).Value = "Sanluqueno"
$ws.Cells.Item(102, 9).Value = 0
$ws.Cells.Item(102, 10).Value = 2.43
$ws.Cells.Item(102, 11).Value = "02/11/2023 08:13"
$ws.Cells.Item(102, 12).Value = 2.17
$ws.Cells.Item(102, 13).Value = "04/11/2023 19:54"
$ws.Cells.Item(102, 14).Value = 3.13
$ws.Cells.Item(102, 15).Value = "02/11/2023 08:13"
$ws.Cells.Item(102, 16).Value = 3.15
$ws.Cells.Item(102, 17).Value = "04/11/2023 19:54"
$ws.Cells.Item(102, 18).Value = 2.77
$ws.Cells.Item(102, 19).Value = "02/11/2023 08:13"
$ws.Cells.Item(102, 20).Value = 3.6
$ws.Cells.Item(102, 21).Value = "04/11/2023 19:54"
$ws.Cells.Item(102, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/cf-intercity-sanluqueno/QVVulhOh/"

# Row 123: restore correct fixture data (F:V)
$ws.Cells.Item(123, 6).Value = "Ceuta"
$ws.Cells.Item(123, 7).Value = 1
$ws.Cells.Item(123, 8).Value = "Granada CF B"
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 1.5
$ws.Cells.Item(123, 11).Value = "17/11/2023 02:12"
$ws.Cells.Item(123, 12).Value = 1.52
$ws.Cells.Item(123, 13).Value = "19/11/2023 11:09"
$ws.Cells.Item(123, 14).Value = 3.74
$ws.Cells.Item(123, 15).Value = "17/11/2023 02:12"
$ws.Cells.Item(123, 16).Value = 3.85
$ws.Cells.Item(123, 17).Value = "19/11/2023 11:09"
$ws.Cells.Item(123, 18).Value = 6.0
$ws.Cells.Item(123, 19).Value = "17/11/2023 02:12"
$ws.Cells.Item(123, 20).Value = 7.18
$ws.Cells.Item(123, 21).Value = "19/11/2023 11:09"
$ws.Cells.Item(123, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/ad-ceuta-granada-cf/bH5wgv6k/"

# Row 124: restore correct fixture data (F:V)
$ws.Cells.Item(124, 6).Value = "Castellon"
$ws.Cells.Item(124, 7).Value = 4
$ws.Cells.Item(124, 8).Value = "Baleares"
$ws.Cells.Item(124, 9).Value = 2
$ws.Cells.Item(124, 10).Value = 1.34
$ws.Cells.Item(124, 11).Value = "16/11/2023 09:13"
$ws.Cells.Item(124, 12).Value = 1.17
$ws.Cells.Item(124, 13).Value = "19/11/2023 11:25"
$ws.Cells.Item(124, 14).Value = 4.47
$ws.Cells.Item(124, 15).Value = "16/11/2023 09:13"
$ws.Cells.Item(124, 16).Value = 6.72
$ws.Cells.Item(124, 17).Value = "19/11/2023 11:25"
$ws.Cells.Item(124, 18).Value = 7.58
$ws.Cells.Item(124, 19).Value = "16/11/2023 09:13"
$ws.Cells.Item(124, 20).Value = 18.15
$ws.Cells.Item(124, 21).Value = "19/11/2023 11:25"
$ws.Cells.Item(124, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/castellon-baleares/nwnYELS8/"

# Row 126: restore correct fixture data (F:V)
$ws.Cells.Item(126, 6).Value = "Real Madrid B"
$ws.Cells.Item(126, 7).Value = 1
$ws.Cells.Item(126, 8).Value = "Algeciras"
$ws.Cells.Item(126, 9).Value = 2
$ws.Cells.Item(126, 10).Value = 2.3
$ws.Cells.Item(126, 11).Value = "19/11/2023 10:27"
$ws.Cells.Item(126, 12).Value = 2.2
$ws.Cells.Item(126, 13).Value = "19/11/2023 19:57"
$ws.Cells.Item(126, 14).Value = 2.94
$ws.Cells.Item(126, 15).Value = "19/11/2023 10:27"
$ws.Cells.Item(126, 16).Value = 3.13
$ws.Cells.Item(126, 17).Value = "19/11/2023 19:57"
$ws.Cells.Item(126, 18).Value = 3.39
$ws.Cells.Item(126, 19).Value = "19/11/2023 10:27"
$ws.Cells.Item(126, 20).Value = 3.54
$ws.Cells.Item(126, 21).Value = "19/11/2023 19:57"
$ws.Cells.Item(126, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/real-madrid-algeciras/8nmUF1D2/"

# Row 127: restore correct fixture data (F:V)
$ws.Cells.Item(127, 6).Value = "Melilla"
$ws.Cells.Item(127, 7).Value = 0
$ws.Cells.Item(127, 8).Value = "UD Ibiza"
$ws.Cells.Item(127, 9).Value = 3
$ws.Cells.Item(127, 10).Value = 4.37
$ws.Cells.Item(127, 11).Value = "16/11/2023 09:13"
$ws.Cells.Item(127, 12).Value = 5.05
$ws.Cells.Item(127, 13).Value = "19/11/2023 19:57"
$ws.Cells.Item(127, 14).Value = 3.2
$ws.Cells.Item(127, 15).Value = "16/11/2023 09:13"
$ws.Cells.Item(127, 16).Value = 3.36
$ws.Cells.Item(127, 17).Value = "19/11/2023 19:58"
$ws.Cells.Item(127, 18).Value = 1.81
$ws.Cells.Item(127, 19).Value = "16/11/2023 09:13"
$ws.Cells.Item(127, 20).Value = 1.78
$ws.Cells.Item(127, 21).Value = "19/11/2023 19:57"
$ws.Cells.Item(127, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/melilla-ud-ibiza/t0oxEurF/"

# Row 128: new fixture appended (Sanluqueno vs Cordoba)
$ws.Cells.Item(127, 1).Copy()
$ws.Cells.Item(128, 1).PasteSpecial(-4122)
$ws.Cells.Item(127, 5).Copy()
$ws.Cells.Item(128, 5).PasteSpecial(-4122)
$ws.Cells.Item(128, 1).Value = 127
$ws.Cells.Item(128, 2).Value = "spain"
$ws.Cells.Item(128, 3).Value = "primera-rfef-group-2"
$ws.Cells.Item(128, 4).Value = "2023-2024"
$ws.Cells.Item(128, 5).Value = 45252.66666666666
$ws.Cells.Item(128, 6).Value = "Sanluqueno"
$ws.Cells.Item(128, 7).Value = 1
$ws.Cells.Item(128, 8).Value = "Cordoba"
$ws.Cells.Item(128, 9).Value = 2
$ws.Cells.Item(128, 10).Value = 2.42
$ws.Cells.Item(128, 11).Value = "19/10/2023 09:13"
$ws.Cells.Item(128, 12).Value = 3.54
$ws.Cells.Item(128, 13).Value = "22/11/2023 15:45"
$ws.Cells.Item(128, 14).Value = 2.91
$ws.Cells.Item(128, 15).Value = "19/10/2023 09:13"
$ws.Cells.Item(128, 16).Value = 3.12
$ws.Cells.Item(128, 17).Value = "22/11/2023 15:49"
$ws.Cells.Item(128, 18).Value = 2.9
$ws.Cells.Item(128, 19).Value = "19/10/2023 09:13"
$ws.Cells.Item(128, 20).Value = 2.2
$ws.Cells.Item(128, 21).Value = "22/11/2023 15:45"
$ws.Cells.Item(128, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/sanluqueno-cordoba/SzSNsH1J/"
